$d = $word.ActiveDocument

# Helper: find literal text range spanning from the start of $startNeedle
# to the end of $endNeedle, then replace its whole .Text (collapses any
# intervening runs into one run that inherits the formatting of the
# first character of the range - matching Word COM Range.Text semantics).
function Replace-Span($startNeedle, $endNeedle, $newText) {
    $rStart = $d.Content
    $rStart.Find.Execute($startNeedle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $startPos = $rStart.Start

    $rEnd = $d.Content
    $rEnd.Find.Execute($endNeedle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $endPos = $rEnd.End

    $target = $d.Range($startPos, $endPos)
    $target.Text = $newText
}

# 1) "4-quintal falconet, otherwise bases, is 8 pans long"
#    -> "Falconet of 4 quintals, alternatively passevolant, is 8 pans long"
$d.Content.Find.Execute("4-quintal falconet, otherwise bases, is 8 pans long", $true, $false, $false, $false, $false, $true, 1, $false, "Falconet of 4 quintals, alternatively passevolant, is 8 pans long", 2) | Out-Null

# 2) "3-quintal <sup>falconet</sup> is 9 pans long" -> "3 quintals, is 9 pans long"
Replace-Span "3-quintal " " is 9 pans long" "3 quintals, is 9 pans long"

# 3) "2-quintal <sup>falconet</sup> is one " -> "2 quintals, is one "
Replace-Span "2-quintal " " is one " "2 quintals, is one "

# 4) "1-quintal is six " -> "1 quintal, is six "
$d.Content.Find.Execute("1-quintal is six ", $true, $false, $false, $false, $false, $true, 1, $false, "1 quintal, is six ", 2) | Out-Null

# 5) "1-quintal and a half is 7" -> "1 quintal and a half is 7"
$d.Content.Find.Execute("1-quintal and a half is 7", $true, $false, $false, $false, $false, $true, 1, $false, "1 quintal and a half is 7", 2) | Out-Null

# 6) "A 60 lb " -> "60 lb "
$d.Content.Find.Execute("A 60 lb ", $true, $false, $false, $false, $false, $true, 1, $false, "60 lb ", 2) | Out-Null

# 7) "</comment>" -> "</comment>," (the musket-comment closing tag run)
$d.Content.Find.Execute("</comment>", $true, $false, $false, $false, $false, $true, 1, $false, "</comment>,", 2) | Out-Null

# 8) "5-quintal is eleven " -> "5 quintals is eleven "
$d.Content.Find.Execute("5-quintal is eleven ", $true, $false, $false, $false, $false, $true, 1, $false, "5 quintals is eleven ", 2) | Out-Null
